# Applies refreshed Universalis market-price data to the Leve profit tables
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value2 = 423.5
$ws.Range("I2").Value2 = 77.59999999999999
$ws.Range("J2").Value2 = 1000
$ws.Range("K2").Value2 = 77.59999999999999
$ws.Range("L2").Value2 = 1000
$ws.Range("M2").Value2 = 35.40000000000001
$ws.Range("N2").Value2 = -1226

# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value2 = 318.14285
$ws.Range("I4").Value2 = 318.14285
$ws.Range("K4").Value2 = 318.14285
$ws.Range("M4").Value2 = -204.14285

# Row 29 (Leve Item ID 4575)
$ws.Range("H29").Value2 = 4510
$ws.Range("I29").Value2 = 1953.6666
$ws.Range("J29").Value2 = 6701.143
$ws.Range("K29").Value2 = 5860.9998
$ws.Range("L29").Value2 = 20103.429
$ws.Range("M29").Value2 = -5579.9998
$ws.Range("N29").Value2 = -20665.429

# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value2 = 13702.833
$ws.Range("J32").Value2 = 13643.4
$ws.Range("L32").Value2 = 13643.4
$ws.Range("N32").Value2 = -14295.4

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value2 = 3544.7778
$ws.Range("J43").Value2 = 4157.5713
$ws.Range("L43").Value2 = 4157.5713
$ws.Range("N43").Value2 = -4295.5713

# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value2 = 1846.3823
$ws.Range("I129").Value2 = 1141.1765
$ws.Range("K129").Value2 = 3423.5295
$ws.Range("M129").Value2 = 1576.4705

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value2 = 47622044
$ws.Range("I45").Value2 = 76924480
$ws.Range("J45").Value2 = 5575.75
$ws.Range("K45").Value2 = 76924480
$ws.Range("L45").Value2 = 5575.75
$ws.Range("M45").Value2 = -76924103
$ws.Range("N45").Value2 = -6329.75

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value2 = 166670380
$ws.Range("I74").Value2 = 166670380
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 166670380
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = -166669506
$ws.Range("N74").ClearContents()

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value2 = 166670380
$ws.Range("I77").Value2 = 166670380
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 833351900
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = -833347532
$ws.Range("N77").ClearContents()

# Row 125 (Leve Item ID 34251)
$ws.Range("H125").Value2 = 0
$ws.Range("J125").Value2 = 0
$ws.Range("L125").Value2 = 0
$ws.Range("N125").ClearContents()

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value2 = 2585.3438
$ws.Range("I86").Value2 = 1974.1072
$ws.Range("J86").Value2 = 6864
$ws.Range("K86").Value2 = 1974.1072
$ws.Range("L86").Value2 = 6864
$ws.Range("M86").Value2 = -851.1071999999999
$ws.Range("N86").Value2 = -9110

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value2 = 2585.3438
$ws.Range("I89").Value2 = 1974.1072
$ws.Range("J89").Value2 = 6864
$ws.Range("K89").Value2 = 9870.536
$ws.Range("L89").Value2 = 34320
$ws.Range("M89").Value2 = -4254.536
$ws.Range("N89").Value2 = -45552

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value2 = 801.1111
$ws.Range("I16").Value2 = 830.4286
$ws.Range("J16").Value2 = 698.5
$ws.Range("K16").Value2 = 830.4286
$ws.Range("L16").Value2 = 698.5
$ws.Range("M16").Value2 = -543.4286
$ws.Range("N16").Value2 = -1272.5

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value2 = 34135.082
$ws.Range("I31").Value2 = 2909.7827
$ws.Range("J31").Value2 = 89379.84
$ws.Range("K31").Value2 = 2909.7827
$ws.Range("L31").Value2 = 89379.84
$ws.Range("M31").Value2 = -2614.7827
$ws.Range("N31").Value2 = -89969.84

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value2 = 34135.082
$ws.Range("I34").Value2 = 2909.7827
$ws.Range("J34").Value2 = 89379.84
$ws.Range("K34").Value2 = 2909.7827
$ws.Range("L34").Value2 = 89379.84
$ws.Range("M34").Value2 = -2707.7827
$ws.Range("N34").Value2 = -89783.84

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value2 = 801.1111
$ws.Range("I113").Value2 = 830.4286
$ws.Range("J113").Value2 = 698.5
$ws.Range("K113").Value2 = 830.4286
$ws.Range("L113").Value2 = 698.5
$ws.Range("M113").Value2 = 1339.5714
$ws.Range("N113").Value2 = -5038.5

# Row 131 (Leve Item ID 35461)
$ws.Range("H131").Value2 = 65000
$ws.Range("I131").Value2 = 0
$ws.Range("J131").Value2 = 65000
$ws.Range("K131").Value2 = 0
$ws.Range("L131").Value2 = 65000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value2 = -75080

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")

# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value2 = 3621255.8
$ws.Range("I4").Value2 = 6666807.5
$ws.Range("J4").Value2 = 1337091.9
$ws.Range("K4").Value2 = 20000422.5
$ws.Range("L4").Value2 = 4011275.7
$ws.Range("M4").Value2 = -20000310.5
$ws.Range("N4").Value2 = -4011499.7

# Row 17 (Leve Item ID 4640)
$ws.Range("H17").Value2 = 435.125
$ws.Range("J17").Value2 = 488.66666
$ws.Range("L17").Value2 = 1465.99998
$ws.Range("N17").Value2 = -1803.99998

# Row 41 (Leve Item ID 4700)
$ws.Range("H41").Value2 = 783
$ws.Range("I41").Value2 = 333
$ws.Range("K41").Value2 = 999
$ws.Range("M41").Value2 = -661

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value2 = 131395
$ws.Range("I132").Value2 = 131395
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 394185
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -391655
$ws.Range("N132").ClearContents()

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 19 (Leve Item ID 2229)
$ws.Range("H19").Value2 = 1051.5
$ws.Range("I19").Value2 = 603
$ws.Range("J19").Value2 = 1500
$ws.Range("K19").Value2 = 603
$ws.Range("L19").Value2 = 1500
$ws.Range("M19").Value2 = -433
$ws.Range("N19").Value2 = -1840

# Row 25 (Leve Item ID 3547)
$ws.Range("H25").Value2 = 6500
$ws.Range("I25").Value2 = 4714.2856
$ws.Range("J25").Value2 = 9000
$ws.Range("K25").Value2 = 4714.2856
$ws.Range("L25").Value2 = 9000
$ws.Range("M25").Value2 = -4484.2856
$ws.Range("N25").Value2 = -9460

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value2 = 8167.3335
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 8167.3335
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 8167.3335
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -8543.333500000001

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value2 = 5357.3887
$ws.Range("I93").Value2 = 1467.0714
$ws.Range("K93").Value2 = 1467.0714
$ws.Range("M93").Value2 = -219.0714

# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value2 = 4177.5
$ws.Range("J100").Value2 = 8400.6
$ws.Range("L100").Value2 = 8400.6
$ws.Range("N100").Value2 = -9482.6

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value2 = 5692.615
$ws.Range("I132").Value2 = 3110.4443
$ws.Range("K132").Value2 = 9331.332900000001
$ws.Range("M132").Value2 = -6801.332900000001

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value2 = 4785.0356
$ws.Range("I136").Value2 = 3088.6
$ws.Range("K136").Value2 = 9265.799999999999
$ws.Range("M136").Value2 = -6715.799999999999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value2 = 3114.611
$ws.Range("I81").Value2 = 1927.7693
$ws.Range("K81").Value2 = 3855.5386
$ws.Range("M81").Value2 = -2794.5386

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value2 = 3114.611
$ws.Range("I84").Value2 = 1927.7693
$ws.Range("K84").Value2 = 19277.693
$ws.Range("M84").Value2 = -13973.693

# Row 135 (Leve Item ID 42043)
$ws.Range("H135").Value2 = 54455.5
$ws.Range("J135").Value2 = 54455.5
$ws.Range("L135").Value2 = 54455.5
$ws.Range("N135").Value2 = -64595.5

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value2 = 4393.579
$ws.Range("I136").Value2 = 2822.7576
$ws.Range("K136").Value2 = 8468.272799999999
$ws.Range("M136").Value2 = -5918.272799999999
